# Generate Report for Handoff
# Updates status text from "Handed back: in sync with en-US" to "Ready for handoff"
# and refreshes the related "Latest ... Datetime" timestamps, then narrows the
# affected columns (their width shrank along with the shorter status text,
# as would happen after an auto-fit in the authoring tool).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-24 00:56:55"

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-24 00:56:51"

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
# de-de's "Latest Handoff Datetime" shares the same underlying value as
# Overview's G2 ("Latest HO Xliff Generate Date"), so it is updated too.
$wsDeDe.Range("H2").Value = "2016-08-24 00:56:55"

# --- Column width adjustments (narrower after shorter status text) ---
$wsOverview.Range("E1").ColumnWidth = 16.38
$wsOverview.Range("F1").ColumnWidth = 16.38
$wsZhCn.Range("C1").ColumnWidth = 16.38
$wsDeDe.Range("C1").ColumnWidth = 16.38
